$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 29.20950566666667
$ws.Range("H2").Value = 87.628517
$ws.Range("I2").Value = 0.01829497698069002
$ws.Range("J2").Value = 0.01840828041918582
$ws.Range("M2").Value = 0.8151449999999999
$ws.Range("N2").Value = 2.445435
$ws.Range("O2").Value = 0.1271069095499719
$ws.Range("P2").Value = 0.1371035811308388
$ws.Range("Q2").Value = 23.809982496655
$ws.Range("R2").Value = 214.289842469895
$ws.Range("S2").Value = 0.002325417984303384
$ws.Range("T2").Value = 0.002523841167931075
$ws.Range("G3").Value = 29.20950566666667
$ws.Range("H3").Value = 87.628517
$ws.Range("I3").Value = 0.01829497698069002
$ws.Range("J3").Value = 0.01840828041918582
$ws.Range("O3").Value = 0.4802730342501803
$ws.Range("P3").Value = 0.5180454245123947
$ws.Range("Q3").Value = 89.96593953546179
$ws.Range("R3").Value = 809.693455819156
$ws.Range("S3").Value = 0.008786584106053198
$ws.Range("T3").Value = 0.009536325444300321
$ws.Range("G4").Value = 29.20950566666667
$ws.Range("H4").Value = 87.628517
$ws.Range("I4").Value = 0.01829497698069002
$ws.Range("J4").Value = 0.01840828041918582
$ws.Range("M4").Value = 0.5185940000000001
$ws.Range("N4").Value = 1.555782
$ws.Range("O4").Value = 0.08086522109705406
$ws.Range("P4").Value = 0.08722508823947427
$ws.Range("Q4").Value = 15.14787438169934
$ws.Range("R4").Value = 136.330869435294
$ws.Range("S4").Value = 0.001479427358509013
$ws.Range("T4").Value = 0.00160566388390047
$ws.Range("G5").Value = 29.20950566666667
$ws.Range("H5").Value = 87.628517
$ws.Range("I5").Value = 0.01829497698069002
$ws.Range("J5").Value = 0.01840828041918582
$ws.Range("M5").Value = 1.402793
$ws.Range("N5").Value = 2.805586
$ws.Range("O5").Value = 0.2187398352051889
$ws.Range("P5").Value = 0.1572954863942594
$ws.Range("Q5").Value = 40.97489008266034
$ws.Range("R5").Value = 245.849340495962
$ws.Range("S5").Value = 0.004001840249838859
$ws.Range("T5").Value = 0.002895539422217754
$ws.Range("G6").Value = 29.20950566666667
$ws.Range("H6").Value = 87.628517
$ws.Range("I6").Value = 0.01829497698069002
$ws.Range("J6").Value = 0.01840828041918582
$ws.Range("M6").Value = 0.5965113333333333
$ws.Range("N6").Value = 1.789534
$ws.Range("O6").Value = 0.09301499989760488
$ws.Range("P6").Value = 0.1003304197230327
$ws.Range("Q6").Value = 17.42380117123089
$ws.Range("R6").Value = 156.814210541078
$ws.Range("S6").Value = 0.001701707281985566
$ws.Range("T6").Value = 0.001846910500836198
$ws.Range("I7").Value = 0.913374480506715
$ws.Range("J7").Value = 0.9190311407684336
$ws.Range("M7").Value = 0.8151449999999999
$ws.Range("N7").Value = 2.445435
$ws.Range("O7").Value = 0.1271069095499719
$ws.Range("P7").Value = 0.1371035811308388
$ws.Range("Q7").Value = 1188.710454061255
$ws.Range("R7").Value = 10698.39408655129
$ws.Range("S7").Value = 0.1160962074790196
$ws.Range("T7").Value = 0.1260024605701123
$ws.Range("I8").Value = 0.913374480506715
$ws.Range("J8").Value = 0.9190311407684336
$ws.Range("O8").Value = 0.4802730342501803
$ws.Range("P8").Value = 0.5180454245123947
$ws.Range("S8").Value = 0.4386691331596422
$ws.Range("T8").Value = 0.4760998774594935
$ws.Range("I9").Value = 0.913374480506715
$ws.Range("J9").Value = 0.9190311407684336
$ws.Range("M9").Value = 0.5185940000000001
$ws.Range("N9").Value = 1.555782
$ws.Range("O9").Value = 0.08086522109705406
$ws.Range("P9").Value = 0.08722508823947427
$ws.Range("Q9").Value = 756.2557694808195
$ws.Range("R9").Value = 6806.301925327375
$ws.Range("S9").Value = 0.07386022931058241
$ws.Range("T9").Value = 0.08016257234835132
$ws.Range("I10").Value = 0.913374480506715
$ws.Range("J10").Value = 0.9190311407684336
$ws.Range("M10").Value = 1.402793
$ws.Range("N10").Value = 2.805586
$ws.Range("O10").Value = 0.2187398352051889
$ws.Range("P10").Value = 0.1572954863942594
$ws.Range("Q10").Value = 2045.6663587263
$ws.Range("R10").Value = 12273.9981523578
$ws.Range("S10").Value = 0.1997913833466639
$ws.Range("T10").Value = 0.1445594502986418
$ws.Range("I11").Value = 0.913374480506715
$ws.Range("J11").Value = 0.9190311407684336
$ws.Range("M11").Value = 0.5965113333333333
$ws.Range("N11").Value = 1.789534
$ws.Range("O11").Value = 0.09301499989760488
$ws.Range("P11").Value = 0.1003304197230327
$ws.Range("Q11").Value = 869.8811351346708
$ws.Range("R11").Value = 7828.930216212038
$ws.Range("S11").Value = 0.084957527210807
$ws.Range("T11").Value = 0.09220678009183453
$ws.Range("G12").Value = 57.98602933333333
$ws.Range("H12").Value = 173.958088
$ws.Range("I12").Value = 0.03631876156896331
$ws.Range("J12").Value = 0.03654368891224535
$ws.Range("M12").Value = 0.8151449999999999
$ws.Range("N12").Value = 2.445435
$ws.Range("O12").Value = 0.1271069095499719
$ws.Range("P12").Value = 0.1371035811308388
$ws.Range("Q12").Value = 47.26702188091999
$ws.Range("R12").Value = 425.4031969282799
$ws.Range("S12").Value = 0.004616365541713214
$ws.Range("T12").Value = 0.005010270617600165
$ws.Range("G13").Value = 57.98602933333333
$ws.Range("H13").Value = 173.958088
$ws.Range("I13").Value = 0.03631876156896331
$ws.Range("J13").Value = 0.03654368891224535
$ws.Range("O13").Value = 0.4802730342501803
$ws.Range("P13").Value = 0.5180454245123947
$ws.Range("Q13").Value = 178.5982846966649
$ws.Range("R13").Value = 1607.384562269984
$ws.Range("S13").Value = 0.01744292181893485
$ws.Range("T13").Value = 0.01893129083579303
$ws.Range("G14").Value = 57.98602933333333
$ws.Range("H14").Value = 173.958088
$ws.Range("I14").Value = 0.03631876156896331
$ws.Range("J14").Value = 0.03654368891224535
$ws.Range("M14").Value = 0.5185940000000001
$ws.Range("N14").Value = 1.555782
$ws.Range("O14").Value = 0.08086522109705406
$ws.Range("P14").Value = 0.08722508823947427
$ws.Range("Q14").Value = 30.07120689609067
$ws.Range("R14").Value = 270.640862064816
$ws.Range("S14").Value = 0.002936924684245408
$ws.Range("T14").Value = 0.003187526489966498
$ws.Range("G15").Value = 57.98602933333333
$ws.Range("H15").Value = 173.958088
$ws.Range("I15").Value = 0.03631876156896331
$ws.Range("J15").Value = 0.03654368891224535
$ws.Range("M15").Value = 1.402793
$ws.Range("N15").Value = 2.805586
$ws.Range("O15").Value = 0.2187398352051889
$ws.Range("P15").Value = 0.1572954863942594
$ws.Range("Q15").Value = 81.34239604659466
$ws.Range("R15").Value = 488.0543762795679
$ws.Range("S15").Value = 0.007944359920451583
$ws.Range("T15").Value = 0.005748157322092134
$ws.Range("G16").Value = 57.98602933333333
$ws.Range("H16").Value = 173.958088
$ws.Range("I16").Value = 0.03631876156896331
$ws.Range("J16").Value = 0.03654368891224535
$ws.Range("M16").Value = 0.5965113333333333
$ws.Range("N16").Value = 1.789534
$ws.Range("O16").Value = 0.09301499989760488
$ws.Range("P16").Value = 0.1003304197230327
$ws.Range("Q16").Value = 34.58932367233244
$ws.Range("R16").Value = 311.3039130509919
$ws.Range("S16").Value = 0.003378189603618258
$ws.Range("T16").Value = 0.003666443646793513
$ws.Range("G17").Value = 29.481085
$ws.Range("H17").Value = 58.96217
$ws.Range("I17").Value = 0.01846507700595112
$ws.Range("J17").Value = 0.01238628926567028
$ws.Range("M17").Value = 0.8151449999999999
$ws.Range("N17").Value = 2.445435
$ws.Range("O17").Value = 0.1271069095499719
$ws.Range("P17").Value = 0.1371035811308388
$ws.Range("Q17").Value = 24.031359032325
$ws.Range("R17").Value = 144.18815419395
$ws.Range("S17").Value = 0.002347038872828695
$ws.Range("T17").Value = 0.001698204615245863
$ws.Range("G18").Value = 29.481085
$ws.Range("H18").Value = 58.96217
$ws.Range("I18").Value = 0.01846507700595112
$ws.Range("J18").Value = 0.01238628926567028
$ws.Range("O18").Value = 0.4802730342501803
$ws.Range("P18").Value = 0.5180454245123947
$ws.Range("Q18").Value = 90.80241003792668
$ws.Range("R18").Value = 544.8144602275601
$ws.Range("S18").Value = 0.008868278561311381
$ws.Range("T18").Value = 0.006416660480767476
$ws.Range("G19").Value = 29.481085
$ws.Range("H19").Value = 58.96217
$ws.Range("I19").Value = 0.01846507700595112
$ws.Range("J19").Value = 0.01238628926567028
$ws.Range("M19").Value = 0.5185940000000001
$ws.Range("N19").Value = 1.555782
$ws.Range("O19").Value = 0.08086522109705406
$ws.Range("P19").Value = 0.08722508823947427
$ws.Range("Q19").Value = 15.28871379449
$ws.Range("R19").Value = 91.73228276694002
$ws.Range("S19").Value = 0.001493182534660367
$ws.Range("T19").Value = 0.001080395174157743
$ws.Range("G20").Value = 29.481085
$ws.Range("H20").Value = 58.96217
$ws.Range("I20").Value = 0.01846507700595112
$ws.Range("J20").Value = 0.01238628926567028
$ws.Range("M20").Value = 1.402793
$ws.Range("N20").Value = 2.805586
$ws.Range("O20").Value = 0.2187398352051889
$ws.Range("P20").Value = 0.1572954863942594
$ws.Range("Q20").Value = 41.355859670405
$ws.Range("R20").Value = 165.42343868162
$ws.Range("S20").Value = 0.004039047901332873
$ws.Range("T20").Value = 0.0019483073946636
$ws.Range("G21").Value = 29.481085
$ws.Range("H21").Value = 58.96217
$ws.Range("I21").Value = 0.01846507700595112
$ws.Range("J21").Value = 0.01238628926567028
$ws.Range("M21").Value = 0.5965113333333333
$ws.Range("N21").Value = 1.789534
$ws.Range("O21").Value = 0.09301499989760488
$ws.Range("P21").Value = 0.1003304197230327
$ws.Range("Q21").Value = 17.58580132146333
$ws.Range("R21").Value = 105.51480792878
$ws.Range("S21").Value = 0.00171752913581781
$ws.Range("T21").Value = 0.001242721600835594
$ws.Range("G22").Value = 21.628479
$ws.Range("H22").Value = 64.885437
$ws.Range("I22").Value = 0.01354670393768061
$ws.Range("J22").Value = 0.01363060063446486
$ws.Range("M22").Value = 0.8151449999999999
$ws.Range("N22").Value = 2.445435
$ws.Range("O22").Value = 0.1271069095499719
$ws.Range("P22").Value = 0.1371035811308388
$ws.Range("Q22").Value = 17.630346514455
$ws.Range("R22").Value = 158.673118630095
$ws.Range("S22").Value = 0.001721879672107017
$ws.Range("T22").Value = 0.001868804159949417
$ws.Range("G23").Value = 21.628479
$ws.Range("H23").Value = 64.885437
$ws.Range("I23").Value = 0.01354670393768061
$ws.Range("J23").Value = 0.01363060063446486
$ws.Range("O23").Value = 0.4802730342501803
$ws.Range("P23").Value = 0.5180454245123947
$ws.Range("Q23").Value = 66.616205565524
$ws.Range("R23").Value = 599.5458500897161
$ws.Range("S23").Value = 0.006506116604238732
$ws.Range("T23").Value = 0.007061270292040266
$ws.Range("G24").Value = 21.628479
$ws.Range("H24").Value = 64.885437
$ws.Range("I24").Value = 0.01354670393768061
$ws.Range("J24").Value = 0.01363060063446486
$ws.Range("M24").Value = 0.5185940000000001
$ws.Range("N24").Value = 1.555782
$ws.Range("O24").Value = 0.08086522109705406
$ws.Range("P24").Value = 0.08722508823947427
$ws.Range("Q24").Value = 11.216399438526
$ws.Range("R24").Value = 100.947594946734
$ws.Range("S24").Value = 0.001095457209056875
$ws.Range("T24").Value = 0.001188930343098232
$ws.Range("G25").Value = 21.628479
$ws.Range("H25").Value = 64.885437
$ws.Range("I25").Value = 0.01354670393768061
$ws.Range("J25").Value = 0.01363060063446486
$ws.Range("M25").Value = 1.402793
$ws.Range("N25").Value = 2.805586
$ws.Range("O25").Value = 0.2187398352051889
$ws.Range("P25").Value = 0.1572954863942594
$ws.Range("Q25").Value = 30.340278941847
$ws.Range("R25").Value = 182.041673651082
$ws.Range("S25").Value = 0.002963203786901741
$ws.Range("T25").Value = 0.002144031956644051
$ws.Range("G26").Value = 21.628479
$ws.Range("H26").Value = 64.885437
$ws.Range("I26").Value = 0.01354670393768061
$ws.Range("J26").Value = 0.01363060063446486
$ws.Range("M26").Value = 0.5965113333333333
$ws.Range("N26").Value = 1.789534
$ws.Range("O26").Value = 0.09301499989760488
$ws.Range("P26").Value = 0.1003304197230327
$ws.Range("Q26").Value = 12.901632846262
$ws.Range("R26").Value = 116.114695616358
$ws.Range("S26").Value = 0.001260046665376246
$ws.Range("T26").Value = 0.001367563882732896
